$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.173.06'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '3.094.72'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.22'
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.51'
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.090.64'
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("E9").Value = '  -1.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.46'
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("E11").Value = '  -2.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.476'
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000246'
$ws.Range("E13").Value = '  -1.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.51'
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.121'
$ws.Range("E15").Value = '  -1.77%  '
$ws.Range("D16").Value = '3.604.35'
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("D17").Value = '67.051.02'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.08'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.66'
$ws.Range("E19").Value = '  +3.12%  '
$ws.Range("B20").Value = 'WrappedEther'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D20").Value = '3.090.78'
$ws.Range("E20").Value = '  -1.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '486.83'
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.78'
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.695'
$ws.Range("E23").Value = '  -3.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.44'
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.96'
$ws.Range("E25").Value = '  -2.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.40'
$ws.Range("E27").Value = '  +3.38%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  -2.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.32'
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.23'
$ws.Range("E32").Value = '  -2.21%  '
$ws.Range("E33").Value = '  -1.82%  '
$ws.Range("D34").Value = '0.0₃0940'
$ws.Range("E34").Value = '  -5.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.75'
$ws.Range("E36").Value = '  -3.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.966'
$ws.Range("E37").Value = '  -2.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '46.44'
$ws.Range("E38").Value = '  -2.89%  '
$ws.Range("E39").Value = '  -4.63%  '
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.305'
$ws.Range("E41").Value = '  -2.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.41'
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D43").Value = '2.787.81'
$ws.Range("E43").Value = '  -2.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '380.25'
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.57'
$ws.Range("E45").Value = '  -8.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0349'
$ws.Range("E46").Value = '  -2.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '135.13'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.79'
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("E50").Value = '  -1.71%  '
$ws.Range("E51").Value = '  -1.96%  '
